$wb = $excel.ActiveWorkbook
$aw = $excel.ActiveWindow
$ws = $wb.Worksheets.Item("BOO")
$ws.Activate()
$r = $aw.ScrollWorkbookTabs(1, 2)
Write-Output $r
